$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# IMO numbers for the new "Maersk vessels not in projects" rows (227-400)
$imoNumbers = @(9260445,9214903,9219795,9245744,9219800,9245770,9198575,9245768,9245756,9198587,9193276,9215177,9215189,9215191,9190743,9190755,9190767,9190779,9190731,9190781,9928188,9398395,9175781,9624275,9398424,9164275,9697026,9313917,9313929,9697014,9840702,9928190,9313905,9313931,9313955,9924182,9924194,9924209,9924211,9924223,9924235,9761023,9356074,9302449,9858711,9871165,9299032,9299044,9502867,9458030,9502910,9502958,9502946,9502972,9694567,9501239,9928205,9235567,9235579,9235543,9739680,9739692,9235555,9713375,9727871,9858735,9193264,9329629,9624287,9343077,9343974,9865881,9294161,9360245,9162215,9624299,9193240,9290476,9858747,9840738,9459424,9893864,9386017,9356165,9385996,9356098,9894648,9893876,9894650,9215907,9402627,9402615,9231470,9215878,9215880,9356127,9385972,9402639,9893852,9894674,9894686,9356139,9215919,9222118,9306172,9840714,9740457,9840726,9726671,9344954,9833369,9411381,9433054,9411379,9344681,9436173,9436185,9401697,9408956,9235531,9411367,9840697,9550357,9459632,9391660,9362401,9858723,9289192,9398436,9786085,9885130,9192442,9220885,9192454,9192466,9885166,9192478,9885154,9885142,9885178,9220897,9251614,9251638,9251626,9120865,9166792,9499022,9332690,9823730,9823742,9332705,9634696,9823766,9823728,9823754,9936379,9936410,9936408,9964481,9969041,9928217,9966128,9928229,9943073,9943085,9943097,9966116,9948748,9944546,9948750,9948774,9948815,9948762,9948786)

# Write "Maersk" first so it becomes shared-string index 265 (before "Maersk
# vessels not in projects" = 266), matching the author's original edit order.
$startRow = 227
for ($i = 0; $i -lt $imoNumbers.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 2).Value = $imoNumbers[$i]
    $ws.Cells.Item($r, 4).Value = "Maersk"
    $ws.Cells.Item($r, 3).Value = "Maersk vessels not in projects"
}

# Move the selection/active view to match the saved workbook state.
$ws.Range("E225").Select()
